# New crime data collected — weekly CompStat report refresh:
# bumps the report header (issue number + week-of dates) and refreshes
# every statistic in the crime-complaints table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (rich-text shared strings): issue number and date range.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# ---------------------------------------------------------------------
# Cells that flip between "text placeholder" and "numeric value" need
# their number format fixed up explicitly so the stored cell type
# matches (Excel otherwise keeps whatever type the cell already has).
# ---------------------------------------------------------------------

# Row 15 / Rape: C15 was the text placeholder "0" -> now numeric 1.
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1

# Row 22 / Transit: C22 was numeric 2 -> now the text placeholder "0".
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# Row 30 / Hate Crimes: C30 numeric 1 -> text "0"; D30 text "0" -> numeric 2;
# E30 text "***.*" -> numeric -100.
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 2

$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100

# ---------------------------------------------------------------------
# Plain numeric refreshes across the table (no type/style changes).
# ---------------------------------------------------------------------
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 26
$ws.Range("K14").Value = 36.842105263157
$ws.Range("L14").Value = 52.941176470588
$ws.Range("M14").Value = 23.809523809523
$ws.Range("N14").Value = -61.194029850746
$ws.Range("E15").Value = -66.666666666666
$ws.Range("F15").Value = 10
$ws.Range("H15").Value = -41.176470588235
$ws.Range("I15").Value = 63
$ws.Range("J15").Value = 69
$ws.Range("K15").Value = -8.695652173913
$ws.Range("L15").Value = 10.526315789473
$ws.Range("M15").Value = 53.658536585365
$ws.Range("N15").Value = -61.585365853658
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 29
$ws.Range("E16").Value = -3.448275862068
$ws.Range("G16").Value = 135
$ws.Range("H16").Value = 2.222222222222
$ws.Range("I16").Value = 541
$ws.Range("J16").Value = 590
$ws.Range("K16").Value = -8.305084745762
$ws.Range("L16").Value = 52.394366197183
$ws.Range("M16").Value = -35.595238095238
$ws.Range("N16").Value = -87.698954070031
$ws.Range("C17").Value = 69
$ws.Range("D17").Value = 68
$ws.Range("E17").Value = 1.470588235294
$ws.Range("G17").Value = 259
$ws.Range("H17").Value = -8.108108108108
$ws.Range("I17").Value = 1045
$ws.Range("J17").Value = 996
$ws.Range("K17").Value = 4.919678714859
$ws.Range("L17").Value = 29.171817058096
$ws.Range("M17").Value = 62.015503875969
$ws.Range("N17").Value = -43.726440495422
$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 36
$ws.Range("E18").Value = -19.444444444444
$ws.Range("F18").Value = 106
$ws.Range("G18").Value = 136
$ws.Range("H18").Value = -22.058823529411
$ws.Range("I18").Value = 568
$ws.Range("J18").Value = 621
$ws.Range("K18").Value = -8.534621578099
$ws.Range("L18").Value = 27.069351230425
$ws.Range("M18").Value = -43.984220907297
$ws.Range("N18").Value = -90.303857972004
$ws.Range("C19").Value = 143
$ws.Range("D19").Value = 129
$ws.Range("E19").Value = 10.852713178294
$ws.Range("F19").Value = 487
$ws.Range("G19").Value = 503
$ws.Range("H19").Value = -3.180914512922
$ws.Range("I19").Value = 2021
$ws.Range("J19").Value = 2159
$ws.Range("K19").Value = -6.391848077813
$ws.Range("L19").Value = 51.726726726726
$ws.Range("M19").Value = 30.303030303030
$ws.Range("N19").Value = -21.514563106796
$ws.Range("C20").Value = 43
$ws.Range("D20").Value = 28
$ws.Range("E20").Value = 53.571428571428
$ws.Range("F20").Value = 161
$ws.Range("G20").Value = 111
$ws.Range("H20").Value = 45.045045045045
$ws.Range("I20").Value = 552
$ws.Range("J20").Value = 518
$ws.Range("K20").Value = 6.563706563706
$ws.Range("L20").Value = 64.776119402985
$ws.Range("M20").Value = -12.241653418124
$ws.Range("N20").Value = -92.707094728497
$ws.Range("C21").Value = 314
$ws.Range("D21").Value = 293
$ws.Range("E21").Value = 7.167235494880
$ws.Range("F21").Value = 1146
$ws.Range("G21").Value = 1164
$ws.Range("H21").Value = -1.546391752577
$ws.Range("I21").Value = 4816
$ws.Range("J21").Value = 4972
$ws.Range("K21").Value = -3.137570394207
$ws.Range("L21").Value = 43.675417661097
$ws.Range("M21").Value = 1.581944737397
$ws.Range("N21").Value = -78.584133760227
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = -43.75
$ws.Range("I22").Value = 51
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = -15
$ws.Range("L22").Value = 59.375
$ws.Range("M22").Value = -37.037037037037
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 36
$ws.Range("H23").Value = 5.882352941176
$ws.Range("I23").Value = 146
$ws.Range("J23").Value = 140
$ws.Range("K23").Value = 4.285714285714
$ws.Range("L23").Value = 5.797101449275
$ws.Range("M23").Value = 89.610389610389
$ws.Range("C24").Value = 358
$ws.Range("D24").Value = 295
$ws.Range("E24").Value = 21.355932203389
$ws.Range("F24").Value = 1274
$ws.Range("G24").Value = 1173
$ws.Range("H24").Value = 8.610400682011
$ws.Range("I24").Value = 4987
$ws.Range("J24").Value = 4533
$ws.Range("K24").Value = 10.015442311934
$ws.Range("L24").Value = 47.022405660377
$ws.Range("M24").Value = 43.552101324122
$ws.Range("C25").Value = 107
$ws.Range("D25").Value = 92
$ws.Range("E25").Value = 16.304347826087
$ws.Range("F25").Value = 429
$ws.Range("G25").Value = 401
$ws.Range("H25").Value = 6.982543640897
$ws.Range("I25").Value = 1690
$ws.Range("J25").Value = 1572
$ws.Range("K25").Value = 7.506361323155
$ws.Range("L25").Value = 26.497005988024
$ws.Range("M25").Value = -10.770855332629
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -71.428571428571
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -29.166666666666
$ws.Range("I26").Value = 93
$ws.Range("J26").Value = 107
$ws.Range("K26").Value = -13.084112149532
$ws.Range("L26").Value = -6.060606060606
$ws.Range("C27").Value = 14
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 41
$ws.Range("G27").Value = 42
$ws.Range("H27").Value = -2.380952380952
$ws.Range("I27").Value = 167
$ws.Range("J27").Value = 183
$ws.Range("K27").Value = -8.743169398907
$ws.Range("L27").Value = 5.696202531645
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 400
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 22
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 51
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = -27.142857142857
$ws.Range("L28").Value = -7.272727272727
$ws.Range("M28").Value = -7.272727272727
$ws.Range("N28").Value = -76.818181818181
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 300
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 12
$ws.Range("H29").Value = -16.666666666666
$ws.Range("I29").Value = 42
$ws.Range("J29").Value = 50
$ws.Range("K29").Value = -16
$ws.Range("L29").Value = -14.285714285714
$ws.Range("M29").Value = -8.695652173913
$ws.Range("N29").Value = -78.571428571428
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 15
$ws.Range("J30").Value = 36
$ws.Range("K30").Value = -58.333333333333
$ws.Range("L30").Value = -6.25
